$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.374703407287598
$ws.Range("B1").Value = 2.714120864868164
$ws.Range("C1").Value = 3.353031635284424
$ws.Range("D1").Value = 3.349186420440674
$ws.Range("E1").Value = 1.929690361022949
